# "Add files via upload" — refresh of the HORAS_SOFTWARE tracking sheet:
#  - student names tidied up (drop the joke surnames) in column A
#  - Joaquin Fernández's NIA + weekly hours filled in (row 11)
#  - Gabriel García's NIA filled in (row 13)
#  - Osmar Ali De La Fuente's NIA + weekly hours filled in (row 15)
#  - selection/scroll position left on A13

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Tidy up the student names in column A (rows 10-16) ---
$ws.Cells.Item(10, 1).Value = "Pedro Ramoneda"
$ws.Cells.Item(11, 1).Value = "Joaquin Fernández"
$ws.Cells.Item(12, 1).Value = "Darío Ferrer"
$ws.Cells.Item(13, 1).Value = "Gabriel García"
$ws.Cells.Item(14, 1).Value = "Alejandro Francés"
$ws.Cells.Item(15, 1).Value = "Osmar Ali De La Fuente"
$ws.Cells.Item(16, 1).Value = "Diego Santolaya"

# --- Joaquin Fernández (row 11): NIA + weekly hours S1..S15 ---
$ws.Cells.Item(11, 2).Value = 715821
$joaquinHours = @(7, 8, 4, 3, 10, 8, 8, 8, 6, 4, 3, 2, 5, 15, 9)
for ($i = 0; $i -lt $joaquinHours.Length; $i++) {
    $ws.Cells.Item(11, 3 + $i).Value = $joaquinHours[$i]
}

# --- Gabriel García (row 13): just the NIA, hours still pending ---
$ws.Cells.Item(13, 2).Value = 723553

# --- Osmar Ali De La Fuente (row 15): NIA + weekly hours S1..S15 ---
$ws.Cells.Item(15, 2).Value = 719152
$osmarHours = @(0, 2, 5, 7, 0, 20, 5, 40, 0, 5, 1, 2, 0, 10, 30)
for ($i = 0; $i -lt $osmarHours.Length; $i++) {
    $ws.Cells.Item(15, 3 + $i).Value = $osmarHours[$i]
}

# Leave the cursor/scroll where the author left it when saving.
$ws.Range("A13").Select()
